# Update slide 13 ("Iteración 2 – Certificados de residencia") body text:
#  - HU5 bullet: split into two runs - shorten the first run's text and
#    append a new trailing run with the rest of the new sentence.
#  - HU6 bullet: simplify wording (drop "por correo").

$p = $ppt.ActivePresentation
$slide = $p.Slides.Item(13)
$shape = $slide.Shapes.Item(2)
$textRange = $shape.TextFrame.TextRange

# HU5 (Alta) paragraph: 5th paragraph of this text body.
$hu5Paragraph = $textRange.Paragraphs(5)
$hu5FirstRun = $hu5Paragraph.Runs(1)
$hu5FirstRun.Text = "HU5 (Alta): Como directiva, quiero validar solicitudes "
$hu5SecondRun = $hu5FirstRun.InsertAfter("y enviar certificados de residencia a los vecinos.")

# HU6 (Media) paragraph: 6th paragraph of this text body.
$hu6Paragraph = $textRange.Paragraphs(6)
$hu6FirstRun = $hu6Paragraph.Runs(1)
$hu6FirstRun.Text = "HU6 (Media): Como vecino, quiero recibir mi certificado aprobado para evitar trámites presenciales."
